$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.447.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.831.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.04%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("E7").Value = '  -0.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4066'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.94%  '

$ws.Range("E9").Value = '  +0.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("E11").Value = '  +0.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.337'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.593'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.835.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.11%  '

$ws.Range("E17").Value = '  +0.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001074'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06612'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.069'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.479.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.14%  '

$ws.Range("E25").Value = '  +1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.478'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.048.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.125'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1092'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.689'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.78%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.658'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07188'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2270'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.268'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02345'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.798'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.39%  '

$ws.Range("E40").Value = '  +2.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6272'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.54%  '

$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("B44").Value = 'WEMIXTOKEN'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.412'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.709'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5857'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.994'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.195'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("E51").Value = '  +0.64%  '
